# Auto-generated script applying scheduled market-price updates to Sheets/Tonberry_Profits
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1666.3334
$ws.Range("J28").Value = 1599.8
$ws.Range("L28").Value = 1599.8
$ws.Range("N28").Value = -2569.8
$ws.Range("H31").Value = 457.66666
$ws.Range("I31").Value = 184.5
$ws.Range("J31").Value = 1004
$ws.Range("K31").Value = 553.5
$ws.Range("L31").Value = 3012
$ws.Range("M31").Value = -323.5
$ws.Range("N31").Value = -3472
$ws.Range("H86").Value = 703.9167
$ws.Range("I86").Value = 666.3333
$ws.Range("K86").Value = 666.3333
$ws.Range("M86").Value = 456.6667
$ws.Range("H89").Value = 703.9167
$ws.Range("I89").Value = 666.3333
$ws.Range("K89").Value = 3331.6665
$ws.Range("M89").Value = 2284.3335
$ws.Range("H106").Value = 1724.5
$ws.Range("I106").Value = 1632.6666
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1632.6666
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -1001.6666
$ws.Range("N106").Value = -3262
$ws.Range("H129").Value = 952.8214
$ws.Range("J129").Value = 910.4091
$ws.Range("L129").Value = 2731.2273
$ws.Range("N129").Value = -12731.2273
$ws.Range("H137").Value = 1690.8572
$ws.Range("I137").Value = 1054.5
$ws.Range("K137").Value = 3163.5
$ws.Range("M137").Value = -613.5
$ws.Range("H138").Value = 1931.0588
$ws.Range("J138").Value = 2499.8572
$ws.Range("L138").Value = 7499.571599999999
$ws.Range("N138").Value = -17779.5716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2404.078
$ws.Range("I32").Value = 1566.935
$ws.Range("J32").Value = 7362.5386
$ws.Range("K32").Value = 1566.935
$ws.Range("L32").Value = 7362.5386
$ws.Range("M32").Value = -1279.935
$ws.Range("N32").Value = -7936.5386
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H45").Value = 3335925.2
$ws.Range("I45").Value = 10002420
$ws.Range("K45").Value = 10002420
$ws.Range("M45").Value = -10002043
$ws.Range("H74").Value = 1184.7693
$ws.Range("I74").Value = 534.7059
$ws.Range("K74").Value = 534.7059
$ws.Range("M74").Value = 339.2941
$ws.Range("H77").Value = 1184.7693
$ws.Range("I77").Value = 534.7059
$ws.Range("K77").Value = 2673.5295
$ws.Range("M77").Value = 1694.4705
$ws.Range("H97").Value = 1657.2
$ws.Range("I97").Value = 1561.3572
$ws.Range("J97").Value = 2999
$ws.Range("K97").Value = 1561.3572
$ws.Range("L97").Value = 2999
$ws.Range("M97").Value = -1065.3572
$ws.Range("N97").Value = -3991
$ws.Range("H104").Value = 38625
$ws.Range("J104").Value = 38625
$ws.Range("L104").Value = 38625
$ws.Range("N104").Value = -45613

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5548.4814
$ws.Range("I134").Value = 5783.7085
$ws.Range("K134").Value = 17351.1255
$ws.Range("M134").Value = -14816.1255

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1482.625
$ws.Range("J31").Value = 1824.4667
$ws.Range("L31").Value = 1824.4667
$ws.Range("N31").Value = -2414.4667
$ws.Range("H34").Value = 1482.625
$ws.Range("J34").Value = 1824.4667
$ws.Range("L34").Value = 1824.4667
$ws.Range("N34").Value = -2228.4667
$ws.Range("H60").Value = 4250
$ws.Range("J60").Value = 14000
$ws.Range("L60").Value = 14000
$ws.Range("N60").Value = -15022

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 509.8889
$ws.Range("I34").Value = 323.75
$ws.Range("J34").Value = 1999
$ws.Range("K34").Value = 971.25
$ws.Range("L34").Value = 5997
$ws.Range("M34").Value = -887.25
$ws.Range("N34").Value = -6165
$ws.Range("H68").Value = 1630.4773
$ws.Range("I68").Value = 791.2857
$ws.Range("J68").Value = 1789.2433
$ws.Range("K68").Value = 2373.8571
$ws.Range("L68").Value = 5367.7299
$ws.Range("M68").Value = -1562.8571
$ws.Range("N68").Value = -6989.7299
$ws.Range("H71").Value = 1630.4773
$ws.Range("I71").Value = 791.2857
$ws.Range("J71").Value = 1789.2433
$ws.Range("K71").Value = 7121.571300000001
$ws.Range("L71").Value = 16103.1897
$ws.Range("M71").Value = -3065.571300000001
$ws.Range("N71").Value = -24215.1897
$ws.Range("H107").Value = 1722.0952
$ws.Range("I107").Value = 1158.8182
$ws.Range("K107").Value = 3476.4546
$ws.Range("M107").Value = -1556.4546
$ws.Range("H134").Value = 4483.385
$ws.Range("I134").Value = 4040.5715
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 12121.7145
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -7051.7145
$ws.Range("N134").Value = -25140

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 23666.334
$ws.Range("J88").Value = 27999.5
$ws.Range("L88").Value = 27999.5
$ws.Range("N88").Value = -28901.5
$ws.Range("H91").Value = 23666.334
$ws.Range("J91").Value = 27999.5
$ws.Range("L91").Value = 27999.5
$ws.Range("N91").Value = -31119.5
$ws.Range("H126").Value = 1573770.5
$ws.Range("I126").Value = 1854714.2
$ws.Range("K126").Value = 5564142.6
$ws.Range("M126").Value = -5561672.6
$ws.Range("H132").Value = 1242857.2
$ws.Range("I132").Value = 1749784.8
$ws.Range("J132").Value = 3701.111
$ws.Range("K132").Value = 5249354.4
$ws.Range("L132").Value = 11103.333
$ws.Range("M132").Value = -5246824.4
$ws.Range("N132").Value = -16163.333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 10581.429
$ws.Range("J43").Value = 10581.429
$ws.Range("L43").Value = 10581.429
$ws.Range("N43").Value = -10967.429
$ws.Range("H61").Value = 3888.6
$ws.Range("J61").Value = 4564.3335
$ws.Range("L61").Value = 4564.3335
$ws.Range("N61").Value = -4968.3335
$ws.Range("H68").Value = 1499.75
$ws.Range("I68").Value = 1499.75
$ws.Range("K68").Value = 1499.75
$ws.Range("M68").Value = -750.75
$ws.Range("H71").Value = 1499.75
$ws.Range("I71").Value = 1499.75
$ws.Range("K71").Value = 7498.75
$ws.Range("M71").Value = -3754.75
$ws.Range("H113").Value = 3888.6
$ws.Range("J113").Value = 4564.3335
$ws.Range("L113").Value = 4564.3335
$ws.Range("N113").Value = -8904.333500000001
$ws.Range("H122").Value = 7064.8
$ws.Range("I122").Value = 5947.4
$ws.Range("J122").Value = 9299.6
$ws.Range("K122").Value = 17842.2
$ws.Range("L122").Value = 27898.8
$ws.Range("M122").Value = -15392.2
$ws.Range("N122").Value = -32798.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17369.166
$ws.Range("J64").Value = 15843
$ws.Range("L64").Value = 15843
$ws.Range("N64").Value = -16339
$ws.Range("H67").Value = 17369.166
$ws.Range("J67").Value = 15843
$ws.Range("L67").Value = 15843
$ws.Range("N67").Value = -17559
$ws.Range("H113").Value = 575.75
$ws.Range("J113").Value = 1035.3334
$ws.Range("L113").Value = 3106.0002
$ws.Range("N113").Value = -7446.0002
$ws.Range("H126").Value = 12960.615
$ws.Range("I126").Value = 16072.875
$ws.Range("K126").Value = 48218.625
$ws.Range("M126").Value = -45748.625
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120
$ws.Range("H136").Value = 10896097
$ws.Range("I136").Value = 17923750
$ws.Range("J136").Value = 3234.45
$ws.Range("K136").Value = 53771250
$ws.Range("L136").Value = 9703.349999999999
$ws.Range("M136").Value = -53768700
$ws.Range("N136").Value = -14803.35
